$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2100
$ws.Range("I131").Value = 2095
$ws.Range("J131").Value = 2102.5
$ws.Range("K131").Value = 6285
$ws.Range("L131").Value = 6307.5
$ws.Range("M131").Value = -1245
$ws.Range("N131").Value = -16387.5

$ws.Range("H132").Value = 5132309
$ws.Range("I132").Value = 6537503.5
$ws.Range("J132").Value = 13386.571
$ws.Range("K132").Value = 19612510.5
$ws.Range("L132").Value = 40159.713
$ws.Range("M132").Value = -19609980.5
$ws.Range("N132").Value = -45219.713

$ws.Range("H137").Value = 1833.3334
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -3450
$ws.Range("N137").Value = -9600

$ws.Range("H141").Value = 895
$ws.Range("I141").Value = 895
$ws.Range("K141").Value = 2685
$ws.Range("M141").Value = 2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 87.8
$ws.Range("I5").Value = 34.75
$ws.Range("K5").Value = 34.75
$ws.Range("M5").Value = 77.25

$ws.Range("H62").Value = 74812
$ws.Range("J62").Value = 74812
$ws.Range("L62").Value = 74812
$ws.Range("N62").Value = -76060

$ws.Range("H65").Value = 74812
$ws.Range("J65").Value = 74812
$ws.Range("L65").Value = 224436
$ws.Range("N65").Value = -230676

$ws.Range("H74").Value = 1385.4
$ws.Range("I74").Value = 1206
$ws.Range("K74").Value = 1206
$ws.Range("M74").Value = -332

$ws.Range("H77").Value = 1385.4
$ws.Range("I77").Value = 1206
$ws.Range("K77").Value = 6030
$ws.Range("M77").Value = -1662

$ws.Range("H81").Value = 75000
$ws.Range("J81").Value = 75000
$ws.Range("L81").Value = 75000
$ws.Range("N81").Value = -76996

$ws.Range("H82").Value = 27451.715
$ws.Range("J82").Value = 27451.715
$ws.Range("L82").Value = 27451.715
$ws.Range("N82").Value = -28173.715

$ws.Range("H84").Value = 75000
$ws.Range("J84").Value = 75000
$ws.Range("L84").Value = 225000
$ws.Range("N84").Value = -234984

$ws.Range("H85").Value = 27451.715
$ws.Range("J85").Value = 27451.715
$ws.Range("L85").Value = 27451.715
$ws.Range("N85").Value = -29947.715

$ws.Range("H87").Value = 75000
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77496

$ws.Range("H90").Value = 75000
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -237480

$ws.Range("H134").Value = 33520
$ws.Range("J134").Value = 33520
$ws.Range("L134").Value = 33520
$ws.Range("N134").Value = -43660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 87.8
$ws.Range("I4").Value = 34.75
$ws.Range("K4").Value = 34.75
$ws.Range("M4").Value = 80.25

$ws.Range("H20").Value = 2488.3333
$ws.Range("I20").Value = 2456.3333
$ws.Range("K20").Value = 2456.3333
$ws.Range("M20").Value = -2209.3333

$ws.Range("H99").Value = 55557056
$ws.Range("I99").Value = 62501500
$ws.Range("K99").Value = 62501500
$ws.Range("M99").Value = -62500002

$ws.Range("H107").Value = 1706.6666
$ws.Range("I107").Value = 1216.091
$ws.Range("J107").Value = 2246.3
$ws.Range("K107").Value = 1216.091
$ws.Range("L107").Value = 2246.3
$ws.Range("M107").Value = 703.9090000000001
$ws.Range("N107").Value = -6086.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 957.6667
$ws.Range("I5").Value = 957.6667
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 957.6667
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -845.6667
$ws.Range("N5").ClearContents()

$ws.Range("H31").Value = 2449
$ws.Range("I31").Value = 2449
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2449
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2154
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 2449
$ws.Range("I34").Value = 2449
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2449
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2247
$ws.Range("N34").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 57739.375
$ws.Range("J135").Value = 57739.375
$ws.Range("L135").Value = 57739.375
$ws.Range("N135").Value = -67879.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 812.4167
$ws.Range("J122").Value = 876.2857
$ws.Range("L122").Value = 7886.571300000001
$ws.Range("N122").Value = -12786.5713

$ws.Range("H131").Value = 20409548
$ws.Range("J131").Value = 1608.4
$ws.Range("L131").Value = 4825.200000000001
$ws.Range("N131").Value = -14905.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3907370.2
$ws.Range("I12").Value = 4355190.5
$ws.Range("J12").Value = 2340000
$ws.Range("K12").Value = 4355190.5
$ws.Range("L12").Value = 2340000
$ws.Range("M12").Value = -4355050.5
$ws.Range("N12").Value = -2340280

$ws.Range("H20").Value = 5002000
$ws.Range("I20").Value = 7500000
$ws.Range("K20").Value = 7500000
$ws.Range("M20").Value = -7499755

$ws.Range("H70").Value = 19570192
$ws.Range("I70").Value = 19235086
$ws.Range("K70").Value = 19235086
$ws.Range("M70").Value = -19234816

$ws.Range("H73").Value = 19570192
$ws.Range("I73").Value = 19235086
$ws.Range("K73").Value = 19235086
$ws.Range("M73").Value = -19234150

$ws.Range("H97").Value = 2842.2
$ws.Range("I97").Value = 1166.6666
$ws.Range("J97").Value = 5355.5
$ws.Range("K97").Value = 1166.6666
$ws.Range("L97").Value = 5355.5
$ws.Range("M97").Value = -670.6666
$ws.Range("N97").Value = -6347.5

$ws.Range("H113").Value = 1964.2142
$ws.Range("I113").Value = 1416.625
$ws.Range("J113").Value = 2694.3333
$ws.Range("K113").Value = 1416.625
$ws.Range("L113").Value = 2694.3333
$ws.Range("M113").Value = 753.375
$ws.Range("N113").Value = -7034.3333

$ws.Range("H121").Value = 2517
$ws.Range("J121").Value = 2517
$ws.Range("L121").Value = 2517
$ws.Range("N121").Value = -6011

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4821.2856
$ws.Range("I46").Value = 1035.5714
$ws.Range("J46").Value = 6714.143
$ws.Range("K46").Value = 1035.5714
$ws.Range("L46").Value = 6714.143
$ws.Range("M46").Value = -847.5714
$ws.Range("N46").Value = -7090.143

$ws.Range("H55").Value = 300.26666
$ws.Range("I55").Value = 183.66667
$ws.Range("J55").Value = 766.6667
$ws.Range("K55").Value = 183.66667
$ws.Range("L55").Value = 766.6667
$ws.Range("M55").Value = -10.66667000000001
$ws.Range("N55").Value = -1112.6667

$ws.Range("H61").Value = 1636.2
$ws.Range("I61").Value = 1482.625
$ws.Range("J61").Value = 1811.7142
$ws.Range("K61").Value = 1482.625
$ws.Range("L61").Value = 1811.7142
$ws.Range("M61").Value = -1280.625
$ws.Range("N61").Value = -2215.7142

$ws.Range("H100").Value = 2125
$ws.Range("I100").Value = 1750
$ws.Range("K100").Value = 1750
$ws.Range("M100").Value = -1209

$ws.Range("H113").Value = 1636.2
$ws.Range("I113").Value = 1482.625
$ws.Range("J113").Value = 1811.7142
$ws.Range("K113").Value = 1482.625
$ws.Range("L113").Value = 1811.7142
$ws.Range("M113").Value = 687.375
$ws.Range("N113").Value = -6151.7142

$ws.Range("H136").Value = 18399.334
$ws.Range("I136").Value = 21479.2
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 64437.60000000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -61887.60000000001
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 520.8
$ws.Range("I17").Value = 520.8
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 520.8
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -348.8
$ws.Range("N17").ClearContents()

$ws.Range("H96").Value = 2097.125
$ws.Range("I96").Value = 1749.25
$ws.Range("J96").Value = 2445
$ws.Range("K96").Value = 1749.25
$ws.Range("L96").Value = 2445
$ws.Range("M96").Value = -376.25
$ws.Range("N96").Value = -5191
